$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 2825
$ws.Range("I34").Value = 2825
$ws.Range("K34").Value = 2825
$ws.Range("M34").Value = -2622

$ws.Range("H36").Value = 2825
$ws.Range("I36").Value = 2825
$ws.Range("K36").Value = 2825
$ws.Range("M36").Value = -2110

$ws.Range("H43").Value = 10841.777
$ws.Range("I43").Value = 9900
$ws.Range("J43").Value = 10959.5
$ws.Range("K43").Value = 9900
$ws.Range("L43").Value = 10959.5
$ws.Range("M43").Value = -9831
$ws.Range("N43").Value = -11097.5

$ws.Range("H53").Value = 3634.6428
$ws.Range("I53").Value = 1195.75
$ws.Range("J53").Value = 4610.2
$ws.Range("K53").Value = 1195.75
$ws.Range("L53").Value = 4610.2
$ws.Range("M53").Value = -558.75
$ws.Range("N53").Value = -5884.2

$ws.Range("H76").Value = 7373.6875
$ws.Range("I76").Value = 6498.3335
$ws.Range("J76").Value = 7898.9
$ws.Range("K76").Value = 6498.3335
$ws.Range("L76").Value = 7898.9
$ws.Range("M76").Value = -6183.3335
$ws.Range("N76").Value = -8528.9

$ws.Range("H79").Value = 7373.6875
$ws.Range("I79").Value = 6498.3335
$ws.Range("J79").Value = 7898.9
$ws.Range("K79").Value = 6498.3335
$ws.Range("L79").Value = 7898.9
$ws.Range("M79").Value = -5406.3335
$ws.Range("N79").Value = -10082.9

$ws.Range("H86").Value = 6667
$ws.Range("I86").Value = 9003
$ws.Range("J86").Value = 6375
$ws.Range("K86").Value = 9003
$ws.Range("L86").Value = 6375
$ws.Range("M86").Value = -7880
$ws.Range("N86").Value = -8621

$ws.Range("H89").Value = 6667
$ws.Range("I89").Value = 9003
$ws.Range("J89").Value = 6375
$ws.Range("K89").Value = 45015
$ws.Range("L89").Value = 31875
$ws.Range("M89").Value = -39399
$ws.Range("N89").Value = -43107

$ws.Range("H132").Value = 691.7045000000001
$ws.Range("I132").Value = 583.7805
$ws.Range("K132").Value = 1751.3415
$ws.Range("M132").Value = 778.6585

$ws.Range("H137").Value = 2667.3425
$ws.Range("J137").Value = 2941.5
$ws.Range("L137").Value = 8824.5
$ws.Range("N137").Value = -13924.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9126.049999999999
$ws.Range("I61").Value = 7393.357
$ws.Range("K61").Value = 7393.357
$ws.Range("M61").Value = -7181.357

$ws.Range("H110").Value = 4784.385
$ws.Range("I110").Value = 4395.4814
$ws.Range("J110").Value = 5659.4165
$ws.Range("K110").Value = 4395.4814
$ws.Range("L110").Value = 5659.4165
$ws.Range("M110").Value = -2350.4814
$ws.Range("N110").Value = -9749.416499999999

$ws.Range("H132").Value = 2649.8333
$ws.Range("I132").Value = 1974.5
$ws.Range("K132").Value = 5923.5
$ws.Range("M132").Value = -3393.5

$ws.Range("H136").Value = 9126.049999999999
$ws.Range("I136").Value = 7393.357
$ws.Range("K136").Value = 22180.071
$ws.Range("M136").Value = -19630.071

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 1
$ws.Range("K22").Value = 1
$ws.Range("M22").Value = 172

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 1029.2
$ws.Range("I25").Value = 1099
$ws.Range("J25").Value = 750
$ws.Range("K25").Value = 1099
$ws.Range("L25").Value = 750
$ws.Range("M25").Value = -925
$ws.Range("N25").Value = -1098

$ws.Range("H35").Value = 2683.8
$ws.Range("I35").Value = 2683.8
$ws.Range("K35").Value = 2683.8
$ws.Range("M35").Value = -2389.8

$ws.Range("H58").Value = 2651.3845
$ws.Range("I58").Value = 1329.762
$ws.Range("K58").Value = 1329.762
$ws.Range("M58").Value = -1126.762

$ws.Range("H122").Value = 7432.091
$ws.Range("I122").Value = 2083.1667
$ws.Range("K122").Value = 6249.500100000001
$ws.Range("M122").Value = -3799.500100000001

$ws.Range("H132").Value = 2782.0527
$ws.Range("J132").Value = 4185.4614
$ws.Range("L132").Value = 12556.3842
$ws.Range("N132").Value = -17616.3842

$ws.Range("H136").Value = 2651.3845
$ws.Range("I136").Value = 1329.762
$ws.Range("K136").Value = 3989.286
$ws.Range("M136").Value = -1439.286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 25647376
$ws.Range("I50").Value = 47620280
$ws.Range("J50").Value = 12317.333
$ws.Range("K50").Value = 142860840
$ws.Range("L50").Value = 36951.999
$ws.Range("M50").Value = -142860359
$ws.Range("N50").Value = -37913.999

$ws.Range("H53").Value = 25647376
$ws.Range("I53").Value = 47620280
$ws.Range("J53").Value = 12317.333
$ws.Range("K53").Value = 142860840
$ws.Range("L53").Value = 36951.999
$ws.Range("M53").Value = -142860359
$ws.Range("N53").Value = -37913.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 337160.66
$ws.Range("I80").Value = 557489.4399999999
$ws.Range("K80").Value = 557489.4399999999
$ws.Range("M80").Value = -556491.4399999999

$ws.Range("H83").Value = 337160.66
$ws.Range("I83").Value = 557489.4399999999
$ws.Range("K83").Value = 2787447.2
$ws.Range("M83").Value = -2782455.2

$ws.Range("H86").Value = 20143
$ws.Range("J86").Value = 20143
$ws.Range("L86").Value = 20143
$ws.Range("N86").Value = -22515

$ws.Range("H89").Value = 20143
$ws.Range("J89").Value = 20143
$ws.Range("L89").Value = 60429
$ws.Range("N89").Value = -72285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1583.05
$ws.Range("I16").Value = 450.7857
$ws.Range("J16").Value = 4225
$ws.Range("K16").Value = 450.7857
$ws.Range("L16").Value = 4225
$ws.Range("M16").Value = -280.7857
$ws.Range("N16").Value = -4565

$ws.Range("H32").Value = 3319.6667
$ws.Range("I32").Value = 3319.6667
$ws.Range("K32").Value = 3319.6667
$ws.Range("M32").Value = -3002.6667

$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31352

$ws.Range("H122").Value = 373392.9
$ws.Range("I122").Value = 508237.38
$ws.Range("J122").Value = 13807.667
$ws.Range("K122").Value = 1524712.14
$ws.Range("L122").Value = 41423.001
$ws.Range("M122").Value = -1522262.14
$ws.Range("N122").Value = -46323.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1874.9111
$ws.Range("I122").Value = 1269.8823
$ws.Range("K122").Value = 3809.6469
$ws.Range("M122").Value = -1359.6469

$ws.Range("H132").Value = 1523.5883
$ws.Range("I132").Value = 1144.8445
$ws.Range("J132").Value = 4364.1665
$ws.Range("K132").Value = 3434.5335
$ws.Range("L132").Value = 13092.4995
$ws.Range("M132").Value = -904.5334999999995
$ws.Range("N132").Value = -18152.4995
